$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 loses its wrap-text formatting - copy A1's (non-wrapping) format onto A2
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null

# A3 switches from a numeric value to a text value
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "14000031"

# DateTime values are reformatted from DD-MM-YYYY to YYYY-MM-DD
$ws.Range("B2").Value = "2015-06-21 17:12:10"
$ws.Range("B3").Value = "2015-06-04 08:12:21"

# Column B width is adjusted slightly
$ws.Columns("B").ColumnWidth = 14.6

# Active selection moves to B1
$ws.Range("B1").Select() | Out-Null
